$wb = $excel.ActiveWorkbook

# --- Update the "Date" metadata value (Metadata!B8) ---
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("B8").Value = "2023-05-30T18:28:01+00:00"

# --- Add a new concept row (TEBA) to the Concepts sheet ---
$ws = $wb.Worksheets.Item("Concepts")

# Copy formatting (style) of the last existing data row onto the new row.
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)  # xlPasteFormats

# Bring over the "1" level value as a text/shared-string cell (matches A2:A6).
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4163)     # xlPasteValues

$ws.Range("B7").Value = "TEBA"
$ws.Range("C7").Value = "Tumoral Exome Bioinformatic Analysis"
$ws.Range("D7").Value = ""

$excel.CutCopyMode = $false
